$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.121.45'

$ws.Range("D3").Value = '1.779.53'
$ws.Range("E3").Value = '  -1.86%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.40%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.29'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.08%  '

$ws.Range("E6").Value = '  +0.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3840'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3420'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.68%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.08'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.189'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07462'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.64'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.422'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.76%  '

$ws.Range("D15").Value = '1.784.65'
$ws.Range("E15").Value = '  -1.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.059'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001085'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06662'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.85%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.41'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.47%  '

$ws.Range("E20").Value = '  +0.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.598'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.28'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.58%  '

$ws.Range("D23").Value = '27.155.39'
$ws.Range("E23").Value = '  -1.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.24'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -6.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.382'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.517'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -5.84%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.468'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.18'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.04'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.15%  '

$ws.Range("D30").Value = '1.987.70'
$ws.Range("E30").Value = '  -1.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.69'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.03%  '

$ws.Range("E32").Value = '  -1.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.016'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08705'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.16'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.75%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.634'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6879'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.385'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2193'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06289'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.726'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02324'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.60%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.234'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.28'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.28%  '

$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.002'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.30%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6437'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.859'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.139'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.78%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.87'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.31%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07117'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.34%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.70'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.53%  '
